$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (252-255), continuing the daily series in columns A-D.
$newRows = @(
    @{ Row = 252; Date = 44326; B = 0; C = 3; D = 131.3485113835376 },
    @{ Row = 253; Date = 44327; B = 0; C = 1; D = 43.78283712784589 },
    @{ Row = 254; Date = 44328; B = 0; C = 1; D = 43.78283712784589 },
    @{ Row = 255; Date = 44329; B = 0; C = 0; D = 0 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $r.Date
    $ws.Cells.Item($row, 2).Value2 = $r.B
    $ws.Cells.Item($row, 3).Value2 = $r.C
    $ws.Cells.Item($row, 4).Value2 = $r.D
}

# Copy the formatting (date style, borders, alignment) from the last
# existing row so the new rows match the rest of the table exactly.
$ws.Range("A251:D251").Copy()
$ws.Range("A252:D255").PasteSpecial(-4122)
$excel.CutCopyMode = $false
